# Applies the cryptos.xlsx price/volume/ranking update described by the commit
# "Updated cryptos list on Mon Jun  3 10:32:58 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matching the workbook's existing
# inline-string cells) without perturbing the cell's style/number-format,
# since some "Price" values (e.g. "1.00", "9.11") would otherwise be
# auto-coerced to numbers by Excel's normal cell-value parsing.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "69.155.22"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "3.818.68"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.25%  "
Set-TextValue "D5" "631.75"
$ws.Range("E5").Value = "  +5.74%  "
Set-TextValue "D6" "165.57"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "3.816.81"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +3.44%  "
Set-TextValue "D13" "0.0000251"
$ws.Range("E13").Value = "  +1.29%  "
Set-TextValue "D14" "36.04"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "4.457.12"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "3.859.98"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "69.116.10"
$ws.Range("E17").Value = "  +2.34%  "
Set-TextValue "D18" "18.01"
$ws.Range("E18").Value = "  -1.28%  "
Set-TextValue "D19" "7.13"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("E20").Value = "  -0.19%  "
Set-TextValue "D21" "466.24"
$ws.Range("E21").Value = "  +1.46%  "
Set-TextValue "D22" "9.70"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  +4.82%  "
Set-TextValue "D25" "83.69"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +3.52%  "
Set-TextValue "D28" "10.06"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "3.965.96"
$ws.Range("E30").Value = "  +1.04%  "
Set-TextValue "D31" "2.71"
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("E32").Value = "  +1.99%  "
Set-TextValue "D33" "7.30"
$ws.Range("E33").Value = "  -1.77%  "
Set-TextValue "D34" "29.24"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "9.11"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  +8.35%  "
Set-TextValue "D39" "3.43"
$ws.Range("E39").Value = "  +6.30%  "
Set-TextValue "D40" "5.92"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +0.09%  "
Set-TextValue "D44" "157.64"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("E46").Value = "  +1.48%  "
Set-TextValue "D47" "46.88"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "42.85"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "8.46"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +3.38%  "
Set-TextValue "D51" "0.000281"
$ws.Range("E51").Value = "  +14.42%  "
